# "Clean File With All Test Passed result"
# - Refresh the sample e-mail addresses used on both sheets (imran5/tisha5 -> imran7/tisha7)
# - Set explicit 12.75pt row heights for the header/first-data rows on the
#   "Login Test Credentials" sheet
# - Restore each sheet's remembered cell selection (and keep the original
#   sheet active/selected)

$wb = $excel.ActiveWorkbook

$wsReg = $wb.Worksheets.Item("User Registration Data")
$wsLogin = $wb.Worksheets.Item("Login Test Credentials")

# Update the e-mail addresses (these are shared strings, so updating the
# text here updates every cell that references the string).
$wsReg.Range("D2").Value = "imran7@sample.com"
$wsReg.Range("D3").Value = "tisha7@sample.com"

$wsLogin.Range("A2").Value = "imran7@sample.com"
$wsLogin.Range("A3").Value = "tisha7@sample.com"

# Set explicit row heights on the "Login Test Credentials" sheet rows 1-2
$wsLogin.Rows.Item(1).RowHeight = 12.75
$wsLogin.Rows.Item(2).RowHeight = 12.75

# Restore the remembered selections on each sheet. Select the secondary
# sheet first so that the primary ("User Registration Data") ends up
# active/selected last, matching the original tab state.
$wsLogin.Range("B6").Select()
$wsReg.Activate()
$wsReg.Range("E8").Select()
